# Update countries & provincias Spain
# Refresh the COVID snapshot data: update the "Datos actualizados" timestamp
# and update a handful of country rows whose figures / ordering changed
# between the 09:55 and 11:12 pulls.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: updated timestamp caption
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 14 de Octubre de 2020 a las 11:12"

# Helper table of row -> [Pais, Casos totales, Nuevos casos, Casos activos,
#                         Recuperados, Casos criticos, Muertes hoy, Muertes]
$rows = @{
    5   = @("India", 7241517, 4435, 6301927, 828945, 0, 28, 110645)
    7   = @("Rusia", 1340409, 14231, 1039705, 277499, 0, 239, 23205)
    21  = @("Filipinas", 346536, 1910, 293860, 46227, 0, 78, 6449)
    22  = @("Indonesia", 344749, 4127, 267851, 64742, 0, 129, 12156)
    25  = @("Alemania", 335679, 0, 281900, 44039, 0, 0, 9740)
    35  = @("Polonia", 141804, 6526, 83847, 54740, 0, 116, 3217)
    36  = @("Bolivia", 138922, 227, 102083, 28488, 0, 25, 8351)
    62  = @("Austria", 58672, 1346, 45846, 11954, 0, 11, 872)
    63  = @("Armenia", 58624, 1058, 46713, 10872, 0, 7, 1039)
    64  = @("Singapur", 57889, 5, 57740, 121, 0, 0, 28)
    88  = @("Croacia", 21741, 748, 18197, 3210, 0, 4, 334)
    89  = @("Camerun", 21203, 0, 20117, 663, 0, 0, 423)
    90  = @("Republica de Macedonia", 21193, 0, 16397, 3996, 0, 0, 800)
    102 = @("Finlandia", 12703, 204, 8500, 3857, 0, 0, 346)
    104 = @("Guinea", 11188, 54, 10352, 766, 0, 0, 70)
    121 = @("Lituania", 6505, 139, 2903, 3493, 0, 3, 109)
    216 = @("Montserrat", 13, 0, 12, 0, 0, 0, 1)
    217 = @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
